$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7095420360565186
$ws.Range("B1").Value = 3.983672857284546
$ws.Range("C1").Value = 5.693912982940674
$ws.Range("D1").Value = 1.231135487556458
$ws.Range("E1").Value = 0.7111010551452637
